$wb = $excel.ActiveWorkbook

# --- Reorder sheets: "review_info" moves in front of "hotel_info" ---
$hotelSheet = $wb.Worksheets.Item("hotel_info")
$reviewSheet = $wb.Worksheets.Item("review_info")
$reviewSheet.Move($hotelSheet)

# --- Add a new "State" column to hotel_info, right after "Hotel_Name" (before "City") ---
$ws = $wb.Worksheets.Item("hotel_info")
$ws.Columns.Item(3).Insert()
$ws.Range("C1").Value = "State"
$ws.Range("C2").Value = "Louisiana"
